$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos":
# the "Periodo Mora" list (E16:E22) is rewritten in the opposite order
# (previously ascending 1906..1912, now descending 1912..1906), while each
# period keeps its own "Valor Mora" (F column) amount.
$periods = @("1912", "1911", "1910", "1909", "1908", "1907", "1906")
$values  = @(20267, 38000, 38000, 38000, 38000, 38000, 38000)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

$wb.Save()
